# Refresh the cryptocurrency price / volume(1h) snapshot in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (taken from the updated crypto feed).
$updates = [ordered]@{
    'D2' = '43.329.04'
    'D3' = '2.357.57'
    'E3' = '  +2.81%  '
    'D4' = '1.00'
    'E4' = '  -0.04%  '
    'D5' = '303.64'
    'E5' = '  +1.06%  '
    'D6' = '95.98'
    'E6' = '  +0.35%  '
    'E7' = '  +0.24%  '
    'D8' = '1.00'
    'E8' = '  +0.00%  '
    'E9' = '  +1.25%  '
    'D10' = '34.36'
    'E10' = '  -0.58%  '
    'D11' = '0.0788'
    'E11' = '  +0.38%  '
    'D12' = '18.64'
    'E12' = '  -2.59%  '
    'E13' = '  +2.65%  '
    'E14' = '  +0.27%  '
    'D15' = '2.719.32'
    'E15' = '  +2.49%  '
    'D16' = '2.340.06'
    'E16' = '  +1.85%  '
    'E17' = '  +2.53%  '
    'D18' = '43.278.01'
    'E18' = '  +1.79%  '
    'D19' = '12.26'
    'E19' = '  +0.53%  '
    'E20' = '  +3.76%  '
    'E21' = '  +0.54%  '
    'D22' = '68.28'
    'E22' = '  +1.16%  '
    'D23' = '236.49'
    'E23' = '  +0.64%  '
    'D24' = '2.22'
    'E24' = '  -0.04%  '
    'E25' = '  +0.05%  '
    'E26' = '  +0.82%  '
    'D27' = '24.80'
    'E27' = '  +0.99%  '
    'D28' = '2.36'
    'E28' = '  +15.22%  '
    'D29' = '9.20'
    'E29' = '  +2.16%  '
    'D30' = '31.58'
    'E30' = '  -1.59%  '
    'E31' = '  -0.06%  '
    'D32' = '5.02'
    'E32' = '  +1.36%  '
    'D33' = '0.0730'
    'E33' = '  +5.01%  '
    'D34' = '17.39'
    'E34' = '  +0.05%  '
    'D35' = '1.85'
    'E35' = '  +6.19%  '
    'D36' = '4.39'
    'E36' = '  +0.55%  '
    'E37' = '  -0.61%  '
    'E38' = '  +1.34%  '
    'B39' = 'EnergySwap'
    'C39' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D39' = '22.74'
    'E39' = '  +17.03%  '
    'B40' = 'LidoDAOToken'
    'C40' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D40' = '2.78'
    'E40' = '  +2.15%  '
    'E41' = '  +0.06%  '
    'D42' = '115.17'
    'E42' = '  -29.96%  '
    'D43' = '1.944.83'
    'E43' = '  -0.57%  '
    'D44' = '0.0283'
    'E44' = '  +2.00%  '
    'D45' = '10.01'
    'E45' = '  -4.45%  '
    'E46' = '  +2.03%  '
    'E47' = '  +0.17%  '
    'D48' = '2.582.91'
    'E48' = '  +2.43%  '
    'D49' = '53.33'
    'E49' = '  +0.94%  '
    'D50' = '2.81'
    'E50' = '  -3.33%  '
    'D51' = '72.27'
    'E51' = '  +1.00%  '
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)
    if ($value.Trim() -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Value looks like a plain number (e.g. "1.00", "18.64").
        # The source data keeps it as text, so force a Text format
        # before assigning it, otherwise Excel would silently turn
        # it into a number (and "1.00" would become "1").
        $cell.NumberFormat = "@"
        $cell.Value = $value
        # Drop back to the default style so no stray direct
        # formatting is left behind on the cell.
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
